$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column retains exact text formatting (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '70.866.09'
$ws.Range("D3").Value = '3.652.78'
$ws.Range("E3").Value = '  +6.37%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '580.34'
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").Value = '176.44'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").Value = '3.642.88'
$ws.Range("E7").Value = '  +6.33%  '
$ws.Range("D8").Value = '0.616'
$ws.Range("E8").Value = '  +2.65%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").Value = '0.198'
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").Value = '6.91'
$ws.Range("E11").Value = '  +25.43%  '
$ws.Range("D12").Value = '0.606'
$ws.Range("E12").Value = '  +3.59%  '
$ws.Range("D13").Value = '48.70'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '0.0000291'
$ws.Range("E14").Value = '  +2.58%  '
$ws.Range("D15").Value = '4.243.52'
$ws.Range("E15").Value = '  +6.51%  '
$ws.Range("D16").Value = '671.97'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").Value = '8.93'
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("D18").Value = '3.656.15'
$ws.Range("E18").Value = '  +6.56%  '
$ws.Range("D19").Value = '70.949.48'
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = '11.44'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '0.933'
$ws.Range("E23").Value = '  +3.78%  '
$ws.Range("D24").Value = '17.18'
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("D25").Value = '100.79'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").Value = '3.93'
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").Value = '  +5.06%  '
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  +4.51%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").Value = '35.00'
$ws.Range("E30").Value = '  +4.29%  '
$ws.Range("D31").Value = '3.37'
$ws.Range("E31").Value = '  +3.14%  '
$ws.Range("E32").Value = '  +3.32%  '
$ws.Range("E33").Value = '  -3.05%  '
$ws.Range("D34").Value = '7.40'
$ws.Range("E34").Value = '  +3.22%  '
$ws.Range("D35").Value = '3.99'
$ws.Range("E35").Value = '  +5.83%  '
$ws.Range("D36").Value = '586.41'
$ws.Range("E36").Value = '  +1.91%  '
$ws.Range("D37").Value = '11.10'
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").Value = '0.107'
$ws.Range("E38").Value = '  +3.77%  '
$ws.Range("D39").Value = '58.47'
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '3.610.93'
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("D42").Value = '0.0458'
$ws.Range("E42").Value = '  +8.73%  '
$ws.Range("E43").Value = '  +2.15%  '
$ws.Range("D44").Value = '0.347'
$ws.Range("E44").Value = '  +4.09%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '34.97'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D46").Value = '0.0₃0748'
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("D47").Value = '2.72'
$ws.Range("E47").Value = '  +1.71%  '
$ws.Range("E48").Value = '  +8.43%  '
$ws.Range("D49").Value = '0.133'
$ws.Range("E49").Value = '  +3.11%  '
$ws.Range("D50").Value = '135.24'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("E51").Value = '  +9.55%  '
